# Generate Report for Archive
# Update localization status from "Ready for handoff" to "In Translation"
# across the Overview, zh-cn, and de-de sheets.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        # Cast to [string] explicitly - some cells hold booleans (e.g. "True"),
        # and a bare "$cell.Value() -eq $oldStatus" comparison would coerce the
        # string $oldStatus into a bool and falsely match them.
        if ([string]$cell.Value() -eq $oldStatus) {
            $cell.Value = $newStatus
        }
    }
}

# The Status text got shorter ("Ready for handoff" -> "In Translation"),
# so the Status column no longer needs to be as wide - narrow it to fit
# the new content on every sheet that has a Status column.
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Columns.Item(5).ColumnWidth = 12.5
$ws1.Columns.Item(6).ColumnWidth = 12.5

$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Columns.Item(3).ColumnWidth = 12.5

$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Columns.Item(3).ColumnWidth = 12.5
